$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 88.71429000000001
$ws.Range("I4").Value = 91.833336
$ws.Range("J4").Value = 70
$ws.Range("K4").Value = 91.833336
$ws.Range("L4").Value = 70
$ws.Range("M4").Value = 22.166664
$ws.Range("N4").Value = -298
$ws.Range("H6").Value = 200258.2
$ws.Range("I6").Value = 200258.2
$ws.Range("K6").Value = 600774.6000000001
$ws.Range("M6").Value = -600662.6000000001
$ws.Range("H33").Value = 1468.9259
$ws.Range("I33").Value = 210.13333
$ws.Range("K33").Value = 210.13333
$ws.Range("M33").Value = 18.86667
$ws.Range("H62").Value = 93753800
$ws.Range("I62").Value = 50003580
$ws.Range("K62").Value = 50003580
$ws.Range("M62").Value = -50002956
$ws.Range("H65").Value = 93753800
$ws.Range("I65").Value = 50003580
$ws.Range("K65").Value = 250017900
$ws.Range("M65").Value = -250014780
$ws.Range("H69").Value = 9993.700000000001
$ws.Range("J69").Value = 9988.5
$ws.Range("L69").Value = 29965.5
$ws.Range("N69").Value = -31713.5
$ws.Range("H72").Value = 9993.700000000001
$ws.Range("J72").Value = 9988.5
$ws.Range("L72").Value = 89896.5
$ws.Range("N72").Value = -98632.5
$ws.Range("H98").Value = 3061213.5
$ws.Range("I98").Value = 3369161.8
$ws.Range("K98").Value = 3369161.8
$ws.Range("M98").Value = -3367663.8
$ws.Range("H106").Value = 1999.25
$ws.Range("I106").Value = 1999.25
$ws.Range("K106").Value = 1999.25
$ws.Range("M106").Value = -1368.25
$ws.Range("H107").Value = 18521082
$ws.Range("I107").Value = 10872149
$ws.Range("J107").Value = 62502444
$ws.Range("K107").Value = 10872149
$ws.Range("L107").Value = 62502444
$ws.Range("M107").Value = -10870229
$ws.Range("N107").Value = -62506284
$ws.Range("H112").Value = 6337781.5
$ws.Range("J112").Value = 7746004.5
$ws.Range("L112").Value = 23238013.5
$ws.Range("N112").Value = -23240229.5
$ws.Range("H116").Value = 3960.36
$ws.Range("I116").Value = 3765.9412
$ws.Range("J116").Value = 4373.5
$ws.Range("K116").Value = 3765.9412
$ws.Range("L116").Value = 4373.5
$ws.Range("M116").Value = -323.9412000000002
$ws.Range("N116").Value = -11257.5
$ws.Range("H122").Value = 3061213.5
$ws.Range("I122").Value = 3369161.8
$ws.Range("K122").Value = 10107485.4
$ws.Range("M122").Value = -10105035.4
$ws.Range("H125").Value = 1592.75
$ws.Range("J125").Value = 1936
$ws.Range("L125").Value = 17424
$ws.Range("N125").Value = -22344
$ws.Range("H137").Value = 3025.4443
$ws.Range("I137").Value = 2416.9375
$ws.Range("J137").Value = 3910.5454
$ws.Range("K137").Value = 7250.8125
$ws.Range("L137").Value = 11731.6362
$ws.Range("M137").Value = -4700.8125
$ws.Range("N137").Value = -16831.6362
$ws.Range("H138").Value = 3404.0613
$ws.Range("J138").Value = 3898.3896
$ws.Range("L138").Value = 11695.1688
$ws.Range("N138").Value = -21975.1688

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 72.2
$ws.Range("I5").Value = 85.25
$ws.Range("J5").Value = 20
$ws.Range("K5").Value = 85.25
$ws.Range("L5").Value = 20
$ws.Range("M5").Value = 26.75
$ws.Range("N5").Value = -244
$ws.Range("H6").Value = 9999.666999999999
$ws.Range("I6").Value = 9999
$ws.Range("K6").Value = 9999
$ws.Range("M6").Value = -9826
$ws.Range("H32").Value = 9420.183999999999
$ws.Range("I32").Value = 7283.08
$ws.Range("K32").Value = 7283.08
$ws.Range("M32").Value = -6996.08
$ws.Range("H45").Value = 1740.1666
$ws.Range("I45").Value = 1522.1666
$ws.Range("J45").Value = 1958.1666
$ws.Range("K45").Value = 1522.1666
$ws.Range("L45").Value = 1958.1666
$ws.Range("M45").Value = -1145.1666
$ws.Range("N45").Value = -2712.1666
$ws.Range("H61").Value = 5123.6665
$ws.Range("I61").Value = 2337.1072
$ws.Range("J61").Value = 14876.625
$ws.Range("K61").Value = 2337.1072
$ws.Range("L61").Value = 14876.625
$ws.Range("M61").Value = -2125.1072
$ws.Range("N61").Value = -15300.625
$ws.Range("H98").Value = 35000
$ws.Range("J98").Value = 35000
$ws.Range("L98").Value = 35000
$ws.Range("N98").Value = -40990
$ws.Range("H136").Value = 5123.6665
$ws.Range("I136").Value = 2337.1072
$ws.Range("J136").Value = 14876.625
$ws.Range("K136").Value = 7011.321599999999
$ws.Range("L136").Value = 44629.875
$ws.Range("M136").Value = -4461.321599999999
$ws.Range("N136").Value = -49729.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 72.2
$ws.Range("I4").Value = 85.25
$ws.Range("J4").Value = 20
$ws.Range("K4").Value = 85.25
$ws.Range("L4").Value = 20
$ws.Range("M4").Value = 29.75
$ws.Range("N4").Value = -250
$ws.Range("H20").Value = 2106.8684
$ws.Range("I20").Value = 1442.9286
$ws.Range("K20").Value = 1442.9286
$ws.Range("M20").Value = -1195.9286
$ws.Range("H80").Value = 867.381
$ws.Range("I80").Value = 1616.25
$ws.Range("J80").Value = 406.53845
$ws.Range("K80").Value = 1616.25
$ws.Range("L80").Value = 406.53845
$ws.Range("M80").Value = -618.25
$ws.Range("N80").Value = -2402.53845
$ws.Range("H83").Value = 867.381
$ws.Range("I83").Value = 1616.25
$ws.Range("J83").Value = 406.53845
$ws.Range("K83").Value = 8081.25
$ws.Range("L83").Value = 2032.69225
$ws.Range("M83").Value = -3089.25
$ws.Range("N83").Value = -12016.69225
$ws.Range("H105").Value = 1529.4193
$ws.Range("J105").Value = 1388.1666
$ws.Range("L105").Value = 1388.1666
$ws.Range("N105").Value = -4882.1666
$ws.Range("H107").Value = 694.7
$ws.Range("I107").Value = 712.58826
$ws.Range("J107").Value = 593.3333
$ws.Range("K107").Value = 712.58826
$ws.Range("L107").Value = 593.3333
$ws.Range("M107").Value = 1207.41174
$ws.Range("N107").Value = -4433.3333
$ws.Range("H134").Value = 3299.0278
$ws.Range("I134").Value = 2734.258
$ws.Range("J134").Value = 6800.6
$ws.Range("K134").Value = 8202.773999999999
$ws.Range("L134").Value = 20401.8
$ws.Range("M134").Value = -5667.773999999999
$ws.Range("N134").Value = -25471.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 190.76471
$ws.Range("I7").Value = 158.81818
$ws.Range("J7").Value = 249.33333
$ws.Range("K7").Value = 158.81818
$ws.Range("L7").Value = 249.33333
$ws.Range("M7").Value = -45.81818000000001
$ws.Range("N7").Value = -475.33333
$ws.Range("H22").Value = 4797.6665
$ws.Range("I22").Value = 3971.75
$ws.Range("J22").Value = 6449.5
$ws.Range("K22").Value = 3971.75
$ws.Range("L22").Value = 6449.5
$ws.Range("M22").Value = -3621.75
$ws.Range("N22").Value = -7149.5
$ws.Range("H31").Value = 2355.3447
$ws.Range("I31").Value = 2279.4583
$ws.Range("K31").Value = 2279.4583
$ws.Range("M31").Value = -1984.4583
$ws.Range("H33").Value = 13252.223
$ws.Range("I33").Value = 15924.429
$ws.Range("J33").Value = 3899.5
$ws.Range("K33").Value = 15924.429
$ws.Range("L33").Value = 3899.5
$ws.Range("M33").Value = -15545.429
$ws.Range("N33").Value = -4657.5
$ws.Range("H34").Value = 2355.3447
$ws.Range("I34").Value = 2279.4583
$ws.Range("K34").Value = 2279.4583
$ws.Range("M34").Value = -2077.4583
$ws.Range("H58").Value = 1844.3
$ws.Range("J58").Value = 2455.75
$ws.Range("L58").Value = 2455.75
$ws.Range("N58").Value = -2861.75
$ws.Range("H99").Value = 21771.54
$ws.Range("I99").Value = 16990.428
$ws.Range("J99").Value = 27349.5
$ws.Range("K99").Value = 16990.428
$ws.Range("L99").Value = 27349.5
$ws.Range("M99").Value = -15492.428
$ws.Range("N99").Value = -30345.5
$ws.Range("H126").Value = 21771.54
$ws.Range("I126").Value = 16990.428
$ws.Range("J126").Value = 27349.5
$ws.Range("K126").Value = 50971.284
$ws.Range("L126").Value = 82048.5
$ws.Range("M126").Value = -48501.284
$ws.Range("N126").Value = -86988.5
$ws.Range("H136").Value = 1844.3
$ws.Range("J136").Value = 2455.75
$ws.Range("L136").Value = 7367.25
$ws.Range("N136").Value = -12467.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 3065.125
$ws.Range("I12").Value = 305.7
$ws.Range("J12").Value = 5036.143
$ws.Range("K12").Value = 917.0999999999999
$ws.Range("L12").Value = 15108.429
$ws.Range("M12").Value = -744.0999999999999
$ws.Range("N12").Value = -15454.429
$ws.Range("H70").Value = 824
$ws.Range("I70").Value = 824
$ws.Range("K70").Value = 2472
$ws.Range("M70").Value = -2157
$ws.Range("H73").Value = 824
$ws.Range("I73").Value = 824
$ws.Range("K73").Value = 2472
$ws.Range("M73").Value = -1380
$ws.Range("H75").Value = 600
$ws.Range("I75").Value = 600
$ws.Range("K75").Value = 1800
$ws.Range("M75").Value = -802
$ws.Range("H78").Value = 600
$ws.Range("I78").Value = 600
$ws.Range("K78").Value = 5400
$ws.Range("M78").Value = -408
$ws.Range("H80").Value = 7563.4287
$ws.Range("J80").Value = 7563.4287
$ws.Range("L80").Value = 22690.2861
$ws.Range("N80").Value = -24562.2861
$ws.Range("H83").Value = 7563.4287
$ws.Range("J83").Value = 7563.4287
$ws.Range("L83").Value = 68070.85830000001
$ws.Range("N83").Value = -77430.85830000001
$ws.Range("H94").Value = 146.75
$ws.Range("J94").Value = 198.5
$ws.Range("L94").Value = 595.5
$ws.Range("N94").Value = -1947.5
$ws.Range("H107").Value = 385.96667
$ws.Range("J107").Value = 409.5
$ws.Range("L107").Value = 1228.5
$ws.Range("N107").Value = -5068.5
$ws.Range("H131").Value = 8334926
$ws.Range("J131").Value = 1758.62
$ws.Range("L131").Value = 5275.86
$ws.Range("N131").Value = -15355.86
$ws.Range("H132").Value = 2501.6667
$ws.Range("J132").Value = 2293.3845
$ws.Range("L132").Value = 20640.4605
$ws.Range("N132").Value = -25700.4605
$ws.Range("H133").Value = 9718.091
$ws.Range("I133").Value = 6299.6665
$ws.Range("J133").Value = 11000
$ws.Range("K133").Value = 18898.9995
$ws.Range("L133").Value = 33000
$ws.Range("M133").Value = -13838.9995
$ws.Range("N133").Value = -43120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 3933
$ws.Range("I4").Value = 3933
$ws.Range("K4").Value = 3933
$ws.Range("M4").Value = -3821
$ws.Range("H29").Value = 7500
$ws.Range("I29").Value = 5000
$ws.Range("K29").Value = 5000
$ws.Range("M29").Value = -4710
$ws.Range("H43").Value = 7813.3687
$ws.Range("J43").Value = 17426.5
$ws.Range("L43").Value = 17426.5
$ws.Range("N43").Value = -17728.5
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H70").Value = 7446.4546
$ws.Range("I70").Value = 7605.8125
$ws.Range("K70").Value = 7605.8125
$ws.Range("M70").Value = -7335.8125
$ws.Range("H73").Value = 7446.4546
$ws.Range("I73").Value = 7605.8125
$ws.Range("K73").Value = 7605.8125
$ws.Range("M73").Value = -6669.8125
$ws.Range("H104").Value = 88704
$ws.Range("J104").Value = 88704
$ws.Range("L104").Value = 88704
$ws.Range("N104").Value = -95692

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1966.4
$ws.Range("I16").Value = 1296.1111
$ws.Range("K16").Value = 1296.1111
$ws.Range("M16").Value = -1126.1111
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H68").Value = 3665.3333
$ws.Range("I68").Value = 4368
$ws.Range("K68").Value = 4368
$ws.Range("M68").Value = -3619
$ws.Range("H71").Value = 3665.3333
$ws.Range("I71").Value = 4368
$ws.Range("K71").Value = 21840
$ws.Range("M71").Value = -18096
$ws.Range("H82").Value = 1709.76
$ws.Range("I82").Value = 1816.3478
$ws.Range("J82").Value = 484
$ws.Range("K82").Value = 1816.3478
$ws.Range("L82").Value = 484
$ws.Range("M82").Value = -1455.3478
$ws.Range("N82").Value = -1206
$ws.Range("H85").Value = 1709.76
$ws.Range("I85").Value = 1816.3478
$ws.Range("J85").Value = 484
$ws.Range("K85").Value = 1816.3478
$ws.Range("L85").Value = 484
$ws.Range("M85").Value = -568.3478
$ws.Range("N85").Value = -2980
$ws.Range("H87").Value = 49000
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 49000
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H122").Value = 5810.2354
$ws.Range("I122").Value = 5859.846
$ws.Range("J122").Value = 5649
$ws.Range("K122").Value = 17579.538
$ws.Range("L122").Value = 16947
$ws.Range("M122").Value = -15129.538
$ws.Range("N122").Value = -21847
$ws.Range("H136").Value = 3590.4375
$ws.Range("I136").Value = 1999.9584
$ws.Range("K136").Value = 5999.8752
$ws.Range("M136").Value = -3449.8752

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2264.5625
$ws.Range("I81").Value = 2148.8667
$ws.Range("J81").Value = 4000
$ws.Range("K81").Value = 4297.7334
$ws.Range("L81").Value = 8000
$ws.Range("M81").Value = -3236.7334
$ws.Range("N81").Value = -10122
$ws.Range("H84").Value = 2264.5625
$ws.Range("I84").Value = 2148.8667
$ws.Range("J84").Value = 4000
$ws.Range("K84").Value = 21488.667
$ws.Range("L84").Value = 40000
$ws.Range("M84").Value = -16184.667
$ws.Range("N84").Value = -50608
$ws.Range("H104").Value = 52474.5
$ws.Range("J104").Value = 52474.5
$ws.Range("L104").Value = 52474.5
$ws.Range("N104").Value = -59462.5
$ws.Range("H113").Value = 1005.63336
$ws.Range("J113").Value = 1814.6666
$ws.Range("L113").Value = 5443.9998
$ws.Range("N113").Value = -9783.9998
$ws.Range("H122").Value = 1810.7097
$ws.Range("J122").Value = 2908.1667
$ws.Range("L122").Value = 8724.500100000001
$ws.Range("N122").Value = -13624.5001
$ws.Range("H126").Value = 1594.7693
$ws.Range("I126").Value = 1521.1818
$ws.Range("K126").Value = 4563.5454
$ws.Range("M126").Value = -2093.5454
$ws.Range("H136").Value = 6499.684
$ws.Range("I136").Value = 6683.1113
$ws.Range("K136").Value = 20049.3339
$ws.Range("M136").Value = -17499.3339

